$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 500
$ws.Range("I18").Value = 500
$ws.Range("K18").Value = 500
$ws.Range("M18").Value = -216

$ws.Range("H53").Value = 235
$ws.Range("I53").Value = 142.38461
$ws.Range("J53").Value = 321
$ws.Range("K53").Value = 142.38461
$ws.Range("L53").Value = 321
$ws.Range("M53").Value = 494.61539
$ws.Range("N53").Value = -1595

$ws.Range("H64").Value = 51140772
$ws.Range("J64").Value = 71433570
$ws.Range("L64").Value = 71433570
$ws.Range("N64").Value = -71434066

$ws.Range("H67").Value = 51140772
$ws.Range("J67").Value = 71433570
$ws.Range("L67").Value = 71433570
$ws.Range("N67").Value = -71435286

$ws.Range("H70").Value = 6250
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 6250
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

$ws.Range("H112").Value = 3021.0364
$ws.Range("I112").Value = 950
$ws.Range("J112").Value = 3099.1887
$ws.Range("K112").Value = 2850
$ws.Range("L112").Value = 9297.5661
$ws.Range("M112").Value = -1742
$ws.Range("N112").Value = -11513.5661

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

$ws.Range("H5").Value = 849
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

$ws.Range("H19").Value = 450
$ws.Range("I19").Value = 450
$ws.Range("K19").Value = 450
$ws.Range("M19").Value = -221

$ws.Range("H32").Value = 33095.484
$ws.Range("I32").Value = 31135.684
$ws.Range("K32").Value = 31135.684
$ws.Range("M32").Value = -30848.684

$ws.Range("H43").Value = 11384.4
$ws.Range("J43").Value = 11384.4
$ws.Range("L43").Value = 11384.4
$ws.Range("N43").Value = -12010.4

$ws.Range("H44").Value = 39966.668
$ws.Range("J44").Value = 39966.668
$ws.Range("L44").Value = 39966.668
$ws.Range("N44").Value = -40942.668

$ws.Range("H74").Value = 12501176
$ws.Range("I74").Value = 14706775
$ws.Range("K74").Value = 14706775
$ws.Range("M74").Value = -14705901

$ws.Range("H77").Value = 12501176
$ws.Range("I77").Value = 14706775
$ws.Range("K77").Value = 73533875
$ws.Range("M77").Value = -73529507

$ws.Range("H97").Value = 586.4706
$ws.Range("I97").Value = 594.8182
$ws.Range("K97").Value = 594.8182
$ws.Range("M97").Value = -98.81820000000005

$ws.Range("H122").Value = 4459.5713
$ws.Range("I122").Value = 3597.7637
$ws.Range("J122").Value = 7619.533
$ws.Range("K122").Value = 10793.2911
$ws.Range("L122").Value = 22858.599
$ws.Range("M122").Value = -8343.2911
$ws.Range("N122").Value = -27758.599

$ws.Range("H132").Value = 16361.56
$ws.Range("I132").Value = 18788.03
$ws.Range("J132").Value = 11205.3125
$ws.Range("K132").Value = 56364.09
$ws.Range("L132").Value = 33615.9375
$ws.Range("M132").Value = -53834.09
$ws.Range("N132").Value = -38675.9375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 849
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

$ws.Range("H105").Value = 3092.1667
$ws.Range("I105").Value = 2926.25
$ws.Range("K105").Value = 2926.25
$ws.Range("M105").Value = -1179.25

$ws.Range("H134").Value = 3721.2
$ws.Range("I134").Value = 3721.2
$ws.Range("K134").Value = 11163.6
$ws.Range("M134").Value = -8628.599999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4121.963
$ws.Range("I31").Value = 2490.2856
$ws.Range("K31").Value = 2490.2856
$ws.Range("M31").Value = -2195.2856

$ws.Range("H34").Value = 4121.963
$ws.Range("I34").Value = 2490.2856
$ws.Range("K34").Value = 2490.2856
$ws.Range("M34").Value = -2288.2856

$ws.Range("H58").Value = 3175.1765
$ws.Range("I58").Value = 3132.6
$ws.Range("K58").Value = 3132.6
$ws.Range("M58").Value = -2929.6

$ws.Range("H122").Value = 3243.3914
$ws.Range("I122").Value = 1869.091
$ws.Range("K122").Value = 5607.272999999999
$ws.Range("M122").Value = -3157.272999999999

$ws.Range("H132").Value = 23811228
$ws.Range("I132").Value = 37038024
$ws.Range("K132").Value = 111114072
$ws.Range("M132").Value = -111111542

$ws.Range("H134").Value = 1557.5769
$ws.Range("I134").Value = 1549.875
$ws.Range("J134").Value = 1650
$ws.Range("K134").Value = 4649.625
$ws.Range("L134").Value = 4950
$ws.Range("M134").Value = -2114.625
$ws.Range("N134").Value = -10020

$ws.Range("H136").Value = 3175.1765
$ws.Range("I136").Value = 3132.6
$ws.Range("K136").Value = 9397.799999999999
$ws.Range("M136").Value = -6847.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 365.42856
$ws.Range("I86").Value = 365.42856
$ws.Range("K86").Value = 1096.28568
$ws.Range("M86").Value = 89.71432000000004

$ws.Range("H89").Value = 365.42856
$ws.Range("I89").Value = 365.42856
$ws.Range("K89").Value = 3288.85704
$ws.Range("M89").Value = 2639.14296

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 429.5
$ws.Range("I2").Value = 171.88889
$ws.Range("J2").Value = 1009.125
$ws.Range("K2").Value = 171.88889
$ws.Range("L2").Value = 1009.125
$ws.Range("M2").Value = -58.88889
$ws.Range("N2").Value = -1235.125

$ws.Range("H102").Value = 11595579
$ws.Range("I102").Value = 14575652
$ws.Range("K102").Value = 14575652
$ws.Range("M102").Value = -14574030

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 939.44446
$ws.Range("J22").Value = 965
$ws.Range("L22").Value = 965
$ws.Range("N22").Value = -1555

$ws.Range("H27").Value = 939.44446
$ws.Range("J27").Value = 965
$ws.Range("L27").Value = 965
$ws.Range("N27").Value = -1179

$ws.Range("H40").Value = 25365664
$ws.Range("I40").Value = 13891664
$ws.Range("K40").Value = 13891664
$ws.Range("M40").Value = -13891528

$ws.Range("H68").Value = 2409.0435
$ws.Range("I68").Value = 2180.6843
$ws.Range("J68").Value = 3493.75
$ws.Range("K68").Value = 2180.6843
$ws.Range("L68").Value = 3493.75
$ws.Range("M68").Value = -1431.6843
$ws.Range("N68").Value = -4991.75

$ws.Range("H71").Value = 2409.0435
$ws.Range("I71").Value = 2180.6843
$ws.Range("J71").Value = 3493.75
$ws.Range("K71").Value = 10903.4215
$ws.Range("L71").Value = 17468.75
$ws.Range("M71").Value = -7159.4215
$ws.Range("N71").Value = -24956.75

$ws.Range("H93").Value = 3954.1428
$ws.Range("I93").Value = 2946.5
$ws.Range("K93").Value = 2946.5
$ws.Range("M93").Value = -1698.5

$ws.Range("H132").Value = 6498.3906
$ws.Range("I132").Value = 6279.6514
$ws.Range("J132").Value = 6946.2856
$ws.Range("K132").Value = 18838.9542
$ws.Range("L132").Value = 20838.8568
$ws.Range("M132").Value = -16308.9542
$ws.Range("N132").Value = -25898.8568

$ws.Range("H136").Value = 5732
$ws.Range("I136").Value = 4638.6
$ws.Range("K136").Value = 13915.8
$ws.Range("M136").Value = -11365.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 9154.625
$ws.Range("I136").Value = 6671.077
$ws.Range("K136").Value = 20013.231
$ws.Range("M136").Value = -17463.231
